$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L16").Value = 9
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 8.5
$ws.Range("O16").Value = 1.36
$ws.Range("P16").Value = 3
$ws.Range("Q16").Value = 2.1
$ws.Range("R16").Value = 1.7
$ws.Range("U16").Value = 2.63
$ws.Range("V16").Value = 1.44
$ws.Range("W16").Value = 5
$ws.Range("AM16").Value = 81
$ws.Range("G17").Value = 1.95
$ws.Range("U17").Value = 2.25
$ws.Range("V17").Value = 1.57
$ws.Range("X17").Value = 8
$ws.Range("Z17").Value = 17
$ws.Range("AC17").Value = 6
$ws.Range("AP17").Value = 29
$ws.Range("G18").Value = 1.7
$ws.Range("I18").Value = 5.5
$ws.Range("J18").Value = 2.38
$ws.Range("M18").Value = 1.07
$ws.Range("N18").Value = 9
$ws.Range("Q18").Value = 2.15
$ws.Range("R18").Value = 1.67
$ws.Range("AA18").Value = 15
$ws.Range("AC18").Value = 8
$ws.Range("AJ18").Value = 19
$ws.Range("AO18").Value = 9
$ws.Range("AQ18").Value = 29
$ws.Range("BA18").Value = 151
$ws.Range("O21").Value = 1.67
$ws.Range("P21").Value = 2.1
$ws.Range("G25").Value = 1.75
$ws.Range("H25").Value = 3.4
$ws.Range("I25").Value = 4.75
$ws.Range("J25").Value = 2.5
$ws.Range("K25").Value = 2.05
$ws.Range("L25").Value = 5.5
$ws.Range("O25").Value = 1.44
$ws.Range("P25").Value = 2.63
$ws.Range("W25").Value = 5.5
$ws.Range("X25").Value = 7.5
$ws.Range("Z25").Value = 13
$ws.Range("AA25").Value = 17
$ws.Range("AH25").Value = 10
$ws.Range("AI25").Value = 23
$ws.Range("AJ25").Value = 17
$ws.Range("AK25").Value = 51
$ws.Range("AN25").Value = 3.6
$ws.Range("AO25").Value = 10
$ws.Range("AP25").Value = 23
$ws.Range("AQ25").Value = 34
$ws.Range("AW25").Value = 6.5
$ws.Range("AX25").Value = 29
$ws.Range("AY25").Value = 41
$ws.Range("AZ25").Value = 101
$ws.Range("BA25").Value = 151
$ws.Range("G48").Value = 2.38
$ws.Range("I48").Value = 2.88
$ws.Range("K48").Value = 2.25
$ws.Range("L48").Value = 3.4
$ws.Range("N48").Value = 12
$ws.Range("U48").Value = 1.62
$ws.Range("V48").Value = 2.2
$ws.Range("W48").Value = 9.5
$ws.Range("Y48").Value = 9.5
$ws.Range("Z48").Value = 23
$ws.Range("AC48").Value = 12
$ws.Range("AG48").Value = 151
$ws.Range("AH48").Value = 11
$ws.Range("AK48").Value = 29
$ws.Range("AP48").Value = 21
$ws.Range("BB48").Value = 151
$ws.Range("BC48").Value = 501
$ws.Range("J49").Value = 3.5
$ws.Range("AB49").Value = 29
$ws.Range("AI49").Value = 12
$ws.Range("G50").Value = 1.42
$ws.Range("H50").Value = 4.75
$ws.Range("I50").Value = 7.5
$ws.Range("J50").Value = 1.91
$ws.Range("L50").Value = 6.5
$ws.Range("N50").Value = 17
$ws.Range("O50").Value = 1.17
$ws.Range("P50").Value = 5
$ws.Range("Q50").Value = 1.57
$ws.Range("R50").Value = 2.35
$ws.Range("S50").Value = 1.29
$ws.Range("T50").Value = 3.5
$ws.Range("X50").Value = 7.5
$ws.Range("Z50").Value = 10
$ws.Range("AB50").Value = 21
$ws.Range("AD50").Value = 9
$ws.Range("AH50").Value = 21
$ws.Range("AI50").Value = 41
$ws.Range("AJ50").Value = 21
$ws.Range("AK50").Value = 81
$ws.Range("AL50").Value = 51
$ws.Range("AN50").Value = 3.5
$ws.Range("AP50").Value = 15
$ws.Range("AQ50").Value = 19
$ws.Range("AT50").Value = 3.5
$ws.Range("AU50").Value = 8.5
$ws.Range("AW50").Value = 8
$ws.Range("AX50").Value = 34
$ws.Range("AZ50").Value = 126
$ws.Range("BA50").Value = 126
$ws.Range("J51").Value = 2.5
$ws.Range("K51").Value = 2.07
$ws.Range("L51").Value = 4.25
$ws.Range("T51").Value = 2.55
$ws.Range("W51").Value = 7.4
$ws.Range("X51").Value = 9.5
$ws.Range("Z51").Value = 17
$ws.Range("AA51").Value = 15
$ws.Range("AB51").Value = 25
$ws.Range("AD51").Value = 6.4
$ws.Range("AE51").Value = 14
$ws.Range("AF51").Value = 65
$ws.Range("AG51").Value = 500
$ws.Range("AH51").Value = 10.75
$ws.Range("AI51").Value = 21
$ws.Range("AM51").Value = 40
$ws.Range("AO51").Value = 9.75
$ws.Range("AP51").Value = 18
$ws.Range("AQ51").Value = 35
$ws.Range("AR51").Value = 65
$ws.Range("AV51").Value = 65
$ws.Range("AX51").Value = 22
$ws.Range("AY51").Value = 28
$ws.Range("AZ51").Value = 120
$ws.Range("BB51").Value = 350
$ws.Range("H52").Value = 3.7
$ws.Range("I52").Value = 4.25
$ws.Range("J52").Value = 2.25
$ws.Range("K52").Value = 2.2
$ws.Range("L52").Value = 4.5
$ws.Range("O52").Value = 1.24
$ws.Range("P52").Value = 3.3
$ws.Range("W52").Value = 7.6
$ws.Range("AC52").Value = 11.5
$ws.Range("AE52").Value = 14.5
$ws.Range("AG52").Value = 450
$ws.Range("AH52").Value = 13
$ws.Range("AI52").Value = 25
$ws.Range("AL52").Value = 40
$ws.Range("AO52").Value = 8.25
$ws.Range("AQ52").Value = 27
$ws.Range("AU52").Value = 7.2
$ws.Range("AW52").Value = 6
$ws.Range("AY52").Value = 28
$ws.Range("BB52").Value = 350
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 3
$ws.Range("I58").Value = 2.38
$ws.Range("J58").Value = 3.75
$ws.Range("K58").Value = 1.95
$ws.Range("L58").Value = 3.25
$ws.Range("O58").Value = 1.44
$ws.Range("P58").Value = 2.63
$ws.Range("Q58").Value = 2.4
$ws.Range("R58").Value = 1.53
$ws.Range("S58").Value = 1.53
$ws.Range("T58").Value = 2.38
$ws.Range("X58").Value = 13
$ws.Range("Y58").Value = 12
$ws.Range("AA58").Value = 29
$ws.Range("AC58").Value = 7
$ws.Range("AG58").Value = 1250
$ws.Range("AH58").Value = 6.5
$ws.Range("AO58").Value = 19
$ws.Range("AS58").Value = 301
$ws.Range("AT58").Value = 2.38
$ws.Range("BB58").Value = 251
$ws.Range("G59").Value = 2.2
$ws.Range("I59").Value = 3.2
$ws.Range("M59").Value = 1.08
$ws.Range("N59").Value = 8
$ws.Range("U59").Value = 1.87
$ws.Range("V59").Value = 1.77
$ws.Range("W59").Value = 7
$ws.Range("AA59").Value = 21
$ws.Range("AN59").Value = 4.33
$ws.Range("G85").Value = 2.12
$ws.Range("H85").Value = 3
$ws.Range("I85").Value = 3.55
$ws.Range("J85").Value = 2.75
$ws.Range("L85").Value = 4
$ws.Range("M85").Value = 1.06
$ws.Range("N85").Value = 9.390000000000001
$ws.Range("O85").Value = 1.39
$ws.Range("P85").Value = 2.57
$ws.Range("Q85").Value = 2.12
$ws.Range("R85").Value = 1.57
$ws.Range("T85").Value = 2.4
$ws.Range("U85").Value = 1.85
$ws.Range("V85").Value = 1.75
$ws.Range("W85").Value = 6.2
$ws.Range("X85").Value = 9.25
$ws.Range("Y85").Value = 8.75
$ws.Range("Z85").Value = 20
$ws.Range("AB85").Value = 35
$ws.Range("AC85").Value = 7.5
$ws.Range("AE85").Value = 15
$ws.Range("AF85").Value = 80
$ws.Range("AH85").Value = 9.25
$ws.Range("AI85").Value = 18.5
$ws.Range("AJ85").Value = 12
$ws.Range("AK85").Value = 55
$ws.Range("AL85").Value = 35
$ws.Range("AM85").Value = 45
$ws.Range("AN85").Value = 3.9
$ws.Range("AO85").Value = 11.25
$ws.Range("AP85").Value = 21
$ws.Range("AQ85").Value = 45
$ws.Range("AR85").Value = 90
$ws.Range("AS85").Value = 300
$ws.Range("AU85").Value = 7
$ws.Range("AW85").Value = 5.3
$ws.Range("AX85").Value = 20
$ws.Range("AZ85").Value = 110
$ws.Range("I86").Value = 3.7
$ws.Range("Z86").Value = 17.5
$ws.Range("AA86").Value = 15.5
$ws.Range("AW86").Value = 5.5
$ws.Range("AY86").Value = 26

Write-Host "Applied 227 cell updates"
